# Applies the cryptos.xlsx data refresh described in the commit diff.
# Numeric-looking price strings (e.g. "0.999", "7.38") are written with a
# leading apostrophe quote-prefix, exactly as if a user typed '0.999 into
# the cell in Excel. That forces the cell to stay plain text (matching the
# original inlineStr cells) instead of being coerced into a floating point
# number, which would otherwise introduce binary rounding artifacts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.204.47"
$ws.Range("E2").Value = "  -2.14%  "
$ws.Range("D3").Value = "3.387.04"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'570.29"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").Value = "'140.66"
$ws.Range("E6").Value = "  -5.82%  "
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "3.387.26"
$ws.Range("E8").Value = "  -1.70%  "
$ws.Range("D9").Value = "'0.473"
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("E10").Value = "  -4.22%  "
$ws.Range("E11").Value = "  -1.19%  "
$ws.Range("D12").Value = "'0.391"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "3.963.98"
$ws.Range("E13").Value = "  -1.72%  "
$ws.Range("D14").Value = "'28.10"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").Value = "'0.123"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "3.391.66"
$ws.Range("E16").Value = "  -1.33%  "
$ws.Range("E17").Value = "  -2.72%  "
$ws.Range("D18").Value = "60.373.50"
$ws.Range("E18").Value = "  -1.98%  "
$ws.Range("D19").Value = "'6.25"
$ws.Range("E19").Value = "  -1.23%  "
$ws.Range("D20").Value = "'14.04"
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("D21").Value = "'9.11"
$ws.Range("E21").Value = "  -4.32%  "
$ws.Range("D22").Value = "'388.83"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("E23").Value = "  -1.72%  "
$ws.Range("D24").Value = "'73.45"
$ws.Range("E24").Value = "  +0.86%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  -4.92%  "
$ws.Range("D27").Value = "3.534.69"
$ws.Range("E27").Value = "  -1.49%  "
$ws.Range("E28").Value = "  -1.01%  "
$ws.Range("D30").Value = "'7.38"
$ws.Range("E30").Value = "  -5.55%  "
$ws.Range("D31").Value = "'8.00"
$ws.Range("E31").Value = "  -3.32%  "
$ws.Range("E32").Value = "  -1.68%  "
$ws.Range("D33").Value = "'1.41"
$ws.Range("E33").Value = "  -7.49%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "'23.71"
$ws.Range("E35").Value = "  -1.19%  "
$ws.Range("D36").Value = "'6.93"
$ws.Range("E36").Value = "  -2.19%  "
$ws.Range("D37").Value = "3.416.55"
$ws.Range("E37").Value = "  -1.51%  "
$ws.Range("D38").Value = "'168.07"
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("E39").Value = "  -6.95%  "
$ws.Range("E40").Value = "  -4.62%  "
$ws.Range("D41").Value = "'0.0776"
$ws.Range("E41").Value = "  -2.32%  "
$ws.Range("D42").Value = "'27.20"
$ws.Range("E42").Value = "  +3.12%  "
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").Value = "'4.44"
$ws.Range("E45").Value = "  -1.42%  "
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("D47").Value = "'41.27"
$ws.Range("E47").Value = "  -2.41%  "
$ws.Range("D48").Value = "2.521.83"
$ws.Range("E48").Value = "  -3.63%  "
$ws.Range("E49").Value = "  -3.98%  "
$ws.Range("D50").Value = "'23.19"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").Value = "'6.82"
$ws.Range("E51").Value = "  -3.98%  "
